# logout() - modified by adding '/' to HOME URL
# Mirrors the authored change: the "Placing COD order" scenario (row 9) now
# also exercises logout() as its final action (new column L), and several
# scenario rows (9, 10, 11) that were previously disabled ("NO") are now
# enabled ("YES") for execution. Row 12's Execution Flag cell is reformatted
# to match the plain "YES" styling used elsewhere. Columns C2:C8 lose their
# redundant explicit font formatting (visually identical, but drops the
# applyFont flag) as a side effect of retyping/clearing those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reset-NormalStyle($rangeAddr) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "General"
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
    $rng.Font.Bold = $false
}

# 1) C2:C8 -- drop the redundant explicit font formatting (same visual font,
#    but matches the "unformatted/default" style used elsewhere in the sheet).
Reset-NormalStyle("C2")
Reset-NormalStyle("C3")
Reset-NormalStyle("C4")
Reset-NormalStyle("C5")
Reset-NormalStyle("C6")
Reset-NormalStyle("C7")
Reset-NormalStyle("C8")

# 2) Flip Execution Flag from NO -> YES for rows 9, 10, 11.
$ws.Range("C9").Value = "YES"
$ws.Range("C10").Value = "YES"
$ws.Range("C11").Value = "YES"

# 3) Row 12's Execution Flag is already "YES" but carried a stray custom
#    font; normalize it to the same style used by the other "YES"/"NO" cells.
Reset-NormalStyle("E12")

# 4) Add the missing logout() step as the final action of the
#    "Placing COD order" scenario (row 9).
$ws.Range("L9").Value = "logout"
